$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 496.69232
$ws.Range("J92").Value = 135
$ws.Range("L92").Value = 135
$ws.Range("N92").Value = -2631
$ws.Range("H98").Value = 8097.0625
$ws.Range("I98").Value = 9382.154
$ws.Range("J98").Value = 2528.3333
$ws.Range("K98").Value = 9382.154
$ws.Range("L98").Value = 2528.3333
$ws.Range("M98").Value = -7884.154
$ws.Range("N98").Value = -5524.3333
$ws.Range("H100").Value = 30123.19
$ws.Range("I100").Value = 38203.5
$ws.Range("K100").Value = 38203.5
$ws.Range("M100").Value = -37662.5
$ws.Range("H111").Value = 980.9091
$ws.Range("I111").Value = 982.2222
$ws.Range("J111").Value = 975
$ws.Range("K111").Value = 2946.6666
$ws.Range("L111").Value = 2925
$ws.Range("M111").Value = 120.3334
$ws.Range("N111").Value = -9059
$ws.Range("H116").Value = 5515.923
$ws.Range("I116").Value = 3709.5715
$ws.Range("J116").Value = 7623.3335
$ws.Range("K116").Value = 3709.5715
$ws.Range("L116").Value = 7623.3335
$ws.Range("M116").Value = -267.5715
$ws.Range("N116").Value = -14507.3335
$ws.Range("H122").Value = 8097.0625
$ws.Range("I122").Value = 9382.154
$ws.Range("J122").Value = 2528.3333
$ws.Range("K122").Value = 28146.462
$ws.Range("L122").Value = 7584.999899999999
$ws.Range("M122").Value = -25696.462
$ws.Range("N122").Value = -12484.9999
$ws.Range("H132").Value = 1260.0238
$ws.Range("I132").Value = 1181.6487
$ws.Range("K132").Value = 3544.9461
$ws.Range("M132").Value = -1014.9461
$ws.Range("H135").Value = 1377.375
$ws.Range("I135").Value = 1377.375
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 12396.375
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -9861.375
$ws.Range("N135").Value = ""
$ws.Range("H137").Value = 16609.893
$ws.Range("I137").Value = 7930.9443
$ws.Range("J137").Value = 24832.053
$ws.Range("K137").Value = 23792.8329
$ws.Range("L137").Value = 74496.159
$ws.Range("M137").Value = -21242.8329
$ws.Range("N137").Value = -79596.159
$ws.Range("H138").Value = 8910.388999999999
$ws.Range("I138").Value = 8297.799999999999
$ws.Range("J138").Value = 9146
$ws.Range("K138").Value = 24893.4
$ws.Range("L138").Value = 27438
$ws.Range("M138").Value = -19753.4
$ws.Range("N138").Value = -37718
$ws.Range("H141").Value = 2579.2144
$ws.Range("I141").Value = 2792.6667
$ws.Range("J141").Value = 1298.5
$ws.Range("K141").Value = 8378.000100000001
$ws.Range("L141").Value = 3895.5
$ws.Range("M141").Value = -3198.000100000001
$ws.Range("N141").Value = -14255.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7276.3076
$ws.Range("I32").Value = 4221.8
$ws.Range("K32").Value = 4221.8
$ws.Range("M32").Value = -3934.8
$ws.Range("H45").Value = 9555.066000000001
$ws.Range("I45").Value = 11829.637
$ws.Range("J45").Value = 3300
$ws.Range("K45").Value = 11829.637
$ws.Range("L45").Value = 3300
$ws.Range("M45").Value = -11452.637
$ws.Range("N45").Value = -4054
$ws.Range("H61").Value = 6505.7427
$ws.Range("I61").Value = 2309.4814
$ws.Range("K61").Value = 2309.4814
$ws.Range("M61").Value = -2097.4814
$ws.Range("H74").Value = 3520.75
$ws.Range("I74").Value = 4387.375
$ws.Range("J74").Value = 1787.5
$ws.Range("K74").Value = 4387.375
$ws.Range("L74").Value = 1787.5
$ws.Range("M74").Value = -3513.375
$ws.Range("N74").Value = -3535.5
$ws.Range("H77").Value = 3520.75
$ws.Range("I77").Value = 4387.375
$ws.Range("J77").Value = 1787.5
$ws.Range("K77").Value = 21936.875
$ws.Range("L77").Value = 8937.5
$ws.Range("M77").Value = -17568.875
$ws.Range("N77").Value = -17673.5
$ws.Range("H97").Value = 1611.15
$ws.Range("I97").Value = 1622.2632
$ws.Range("J97").Value = 1400
$ws.Range("K97").Value = 1622.2632
$ws.Range("L97").Value = 1400
$ws.Range("M97").Value = -1126.2632
$ws.Range("N97").Value = -2392
$ws.Range("H102").Value = 1340.7812
$ws.Range("I102").Value = 1274.4138
$ws.Range("K102").Value = 1274.4138
$ws.Range("M102").Value = 347.5862
$ws.Range("H110").Value = 959.25
$ws.Range("I110").Value = 946
$ws.Range("K110").Value = 946
$ws.Range("M110").Value = 1099
$ws.Range("H122").Value = 904.1667
$ws.Range("I122").Value = 437
$ws.Range("J122").Value = 1371.3334
$ws.Range("K122").Value = 1311
$ws.Range("L122").Value = 4114.0002
$ws.Range("M122").Value = 1139
$ws.Range("N122").Value = -9014.0002
$ws.Range("H132").Value = 5260.694
$ws.Range("I132").Value = 3644.725
$ws.Range("K132").Value = 10934.175
$ws.Range("M132").Value = -8404.174999999999
$ws.Range("H136").Value = 6505.7427
$ws.Range("I136").Value = 2309.4814
$ws.Range("K136").Value = 6928.4442
$ws.Range("M136").Value = -4378.4442

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4669.5
$ws.Range("J20").Value = 4999
$ws.Range("L20").Value = 4999
$ws.Range("N20").Value = -5493
$ws.Range("H22").Value = 99
$ws.Range("I22").Value = 99
$ws.Range("K22").Value = 99
$ws.Range("M22").Value = 74
$ws.Range("H58").Value = 48332.668
$ws.Range("J58").Value = 47499.5
$ws.Range("L58").Value = 47499.5
$ws.Range("N58").Value = -48087.5
$ws.Range("H86").Value = 591627.9
$ws.Range("I86").Value = 1113944
$ws.Range("K86").Value = 1113944
$ws.Range("M86").Value = -1112821
$ws.Range("H89").Value = 591627.9
$ws.Range("I89").Value = 1113944
$ws.Range("K89").Value = 5569720
$ws.Range("M89").Value = -5564104
$ws.Range("H94").Value = 2879.8
$ws.Range("I94").Value = 2850
$ws.Range("J94").Value = 2999
$ws.Range("K94").Value = 2850
$ws.Range("L94").Value = 2999
$ws.Range("M94").Value = -2399
$ws.Range("N94").Value = -3901
$ws.Range("H107").Value = 1717.8889
$ws.Range("I107").Value = 1653.6666
$ws.Range("J107").Value = 1750
$ws.Range("K107").Value = 1653.6666
$ws.Range("L107").Value = 1750
$ws.Range("M107").Value = 266.3334
$ws.Range("N107").Value = -5590
$ws.Range("H134").Value = 13221.114
$ws.Range("I134").Value = 10676.167
$ws.Range("J134").Value = 18773.727
$ws.Range("K134").Value = 32028.501
$ws.Range("L134").Value = 56321.181
$ws.Range("M134").Value = -29493.501
$ws.Range("N134").Value = -61391.181

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 314.27585
$ws.Range("I7").Value = 333.6087
$ws.Range("J7").Value = 240.16667
$ws.Range("K7").Value = 333.6087
$ws.Range("L7").Value = 240.16667
$ws.Range("M7").Value = -220.6087
$ws.Range("N7").Value = -466.16667
$ws.Range("H22").Value = 259.48
$ws.Range("I22").Value = 234.625
$ws.Range("J22").Value = 303.66666
$ws.Range("K22").Value = 234.625
$ws.Range("L22").Value = 303.66666
$ws.Range("M22").Value = 115.375
$ws.Range("N22").Value = -1003.66666
$ws.Range("H31").Value = 68559.664
$ws.Range("I31").Value = 72414.07000000001
$ws.Range("K31").Value = 72414.07000000001
$ws.Range("M31").Value = -72119.07000000001
$ws.Range("H34").Value = 68559.664
$ws.Range("I34").Value = 72414.07000000001
$ws.Range("K34").Value = 72414.07000000001
$ws.Range("M34").Value = -72212.07000000001
$ws.Range("H58").Value = 24003.744
$ws.Range("I58").Value = 26556.25
$ws.Range("J58").Value = 9418
$ws.Range("K58").Value = 26556.25
$ws.Range("L58").Value = 9418
$ws.Range("M58").Value = -26353.25
$ws.Range("N58").Value = -9824
$ws.Range("H62").Value = 378537.5
$ws.Range("I62").Value = 253573.5
$ws.Range("K62").Value = 253573.5
$ws.Range("M62").Value = -252949.5
$ws.Range("H65").Value = 378537.5
$ws.Range("I65").Value = 253573.5
$ws.Range("K65").Value = 1267867.5
$ws.Range("M65").Value = -1264747.5
$ws.Range("H86").Value = 3935.1428
$ws.Range("I86").Value = 2758.8
$ws.Range("J86").Value = 4588.6665
$ws.Range("K86").Value = 2758.8
$ws.Range("L86").Value = 4588.6665
$ws.Range("M86").Value = -1635.8
$ws.Range("N86").Value = -6834.6665
$ws.Range("H89").Value = 3935.1428
$ws.Range("I89").Value = 2758.8
$ws.Range("J89").Value = 4588.6665
$ws.Range("K89").Value = 13794
$ws.Range("L89").Value = 22943.3325
$ws.Range("M89").Value = -8178
$ws.Range("N89").Value = -34175.3325
$ws.Range("H99").Value = 38087.39
$ws.Range("I99").Value = 8586.75
$ws.Range("J99").Value = 41774.97
$ws.Range("K99").Value = 8586.75
$ws.Range("L99").Value = 41774.97
$ws.Range("M99").Value = -7088.75
$ws.Range("N99").Value = -44770.97
$ws.Range("H105").Value = 2567.0454
$ws.Range("I105").Value = 2330.25
$ws.Range("J105").Value = 3198.5
$ws.Range("K105").Value = 2330.25
$ws.Range("L105").Value = 3198.5
$ws.Range("M105").Value = -583.25
$ws.Range("N105").Value = -6692.5
$ws.Range("H107").Value = 750.26086
$ws.Range("I107").Value = 693.9286
$ws.Range("J107").Value = 837.8889
$ws.Range("K107").Value = 693.9286
$ws.Range("L107").Value = 837.8889
$ws.Range("M107").Value = 1226.0714
$ws.Range("N107").Value = -4677.8889
$ws.Range("H126").Value = 38087.39
$ws.Range("I126").Value = 8586.75
$ws.Range("J126").Value = 41774.97
$ws.Range("K126").Value = 25760.25
$ws.Range("L126").Value = 125324.91
$ws.Range("M126").Value = -23290.25
$ws.Range("N126").Value = -130264.91
$ws.Range("H132").Value = 31877.305
$ws.Range("I132").Value = 23500.893
$ws.Range("K132").Value = 70502.679
$ws.Range("M132").Value = -67972.679
$ws.Range("H134").Value = 5428.025
$ws.Range("I134").Value = 2730.9688
$ws.Range("K134").Value = 8192.9064
$ws.Range("M134").Value = -5657.9064
$ws.Range("H136").Value = 24003.744
$ws.Range("I136").Value = 26556.25
$ws.Range("J136").Value = 9418
$ws.Range("K136").Value = 79668.75
$ws.Range("L136").Value = 28254
$ws.Range("M136").Value = -77118.75
$ws.Range("N136").Value = -33354

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 20
$ws.Range("I22").Value = 20
$ws.Range("K22").Value = 60
$ws.Range("M22").Value = 109
$ws.Range("H27").Value = 20
$ws.Range("I27").Value = 20
$ws.Range("K27").Value = 60
$ws.Range("M27").Value = 42
$ws.Range("H34").Value = 4623.75
$ws.Range("I34").Value = 163.33333
$ws.Range("J34").Value = 6110.5557
$ws.Range("K34").Value = 489.99999
$ws.Range("L34").Value = 18331.6671
$ws.Range("M34").Value = -405.99999
$ws.Range("N34").Value = -18499.6671
$ws.Range("H39").Value = 7475
$ws.Range("J39").Value = 8600
$ws.Range("L39").Value = 25800
$ws.Range("N39").Value = -26388
$ws.Range("H55").Value = 3582.8333
$ws.Range("J55").Value = 5057.5713
$ws.Range("L55").Value = 15172.7139
$ws.Range("N55").Value = -15526.7139
$ws.Range("H57").Value = 4098.778
$ws.Range("I57").Value = 944.5
$ws.Range("K57").Value = 2833.5
$ws.Range("M57").Value = -2274.5
$ws.Range("H104").Value = 7196.857
$ws.Range("I104").Value = 5249.5713
$ws.Range("J104").Value = 9144.143
$ws.Range("K104").Value = 15748.7139
$ws.Range("L104").Value = 27432.429
$ws.Range("M104").Value = -13127.7139
$ws.Range("N104").Value = -32674.429
$ws.Range("H107").Value = 1546.091
$ws.Range("I107").Value = 1484.2
$ws.Range("K107").Value = 4452.6
$ws.Range("M107").Value = -2532.6
$ws.Range("H122").Value = 33333484
$ws.Range("I122").Value = 150
$ws.Range("J122").Value = 50000148
$ws.Range("K122").Value = 1350
$ws.Range("L122").Value = 450001332
$ws.Range("M122").Value = 1100
$ws.Range("N122").Value = -450006232
$ws.Range("H128").Value = 499978.25
$ws.Range("I128").Value = 499978.25
$ws.Range("K128").Value = 1499934.75
$ws.Range("M128").Value = -1494954.75
$ws.Range("H129").Value = 2634.875
$ws.Range("I129").Value = 1264
$ws.Range("K129").Value = 3792
$ws.Range("M129").Value = 1208
$ws.Range("H132").Value = 7693902
$ws.Range("I132").Value = 1616
$ws.Range("K132").Value = 14544
$ws.Range("M132").Value = -12014

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11347.429
$ws.Range("I70").Value = 11163.167
$ws.Range("K70").Value = 11163.167
$ws.Range("M70").Value = -10893.167
$ws.Range("H73").Value = 11347.429
$ws.Range("I73").Value = 11163.167
$ws.Range("K73").Value = 11163.167
$ws.Range("M73").Value = -10227.167
$ws.Range("H102").Value = 2774.6667
$ws.Range("I102").Value = 2603.4546
$ws.Range("K102").Value = 2603.4546
$ws.Range("M102").Value = -981.4546
$ws.Range("H126").Value = 3233.4375
$ws.Range("I126").Value = 3182.4
$ws.Range("J126").Value = 3999
$ws.Range("K126").Value = 9547.200000000001
$ws.Range("L126").Value = 11997
$ws.Range("M126").Value = -7077.200000000001
$ws.Range("N126").Value = -16937
$ws.Range("H132").Value = 10119.066
$ws.Range("I132").Value = 11716.409
$ws.Range("J132").Value = 5726.375
$ws.Range("K132").Value = 35149.227
$ws.Range("L132").Value = 17179.125
$ws.Range("M132").Value = -32619.227
$ws.Range("N132").Value = -22239.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1148099.2
$ws.Range("I7").Value = 1338698.6
$ws.Range("K7").Value = 1338698.6
$ws.Range("M7").Value = -1338586.6
$ws.Range("H22").Value = 2062.25
$ws.Range("I22").Value = 2071.2856
$ws.Range("J22").Value = 1999
$ws.Range("K22").Value = 2071.2856
$ws.Range("L22").Value = 1999
$ws.Range("M22").Value = -1776.2856
$ws.Range("N22").Value = -2589
$ws.Range("H27").Value = 2062.25
$ws.Range("I27").Value = 2071.2856
$ws.Range("J27").Value = 1999
$ws.Range("K27").Value = 2071.2856
$ws.Range("L27").Value = 1999
$ws.Range("M27").Value = -1964.2856
$ws.Range("N27").Value = -2213
$ws.Range("H40").Value = 4610.1113
$ws.Range("I40").Value = 4499.7144
$ws.Range("J40").Value = 4996.5
$ws.Range("K40").Value = 4499.7144
$ws.Range("L40").Value = 4996.5
$ws.Range("M40").Value = -4363.7144
$ws.Range("N40").Value = -5268.5
$ws.Range("H55").Value = 155.94118
$ws.Range("I55").Value = 151.4
$ws.Range("J55").Value = 190
$ws.Range("K55").Value = 151.4
$ws.Range("L55").Value = 190
$ws.Range("M55").Value = 21.59999999999999
$ws.Range("N55").Value = -536
$ws.Range("H93").Value = 4898
$ws.Range("I93").Value = 4487.2
$ws.Range("J93").Value = 6952
$ws.Range("K93").Value = 4487.2
$ws.Range("L93").Value = 6952
$ws.Range("M93").Value = -3239.2
$ws.Range("N93").Value = -9448
$ws.Range("H126").Value = 1148099.2
$ws.Range("I126").Value = 1338698.6
$ws.Range("K126").Value = 4016095.8
$ws.Range("M126").Value = -4013625.8
$ws.Range("H132").Value = 7722.478
$ws.Range("I132").Value = 7202.0586
$ws.Range("J132").Value = 9197
$ws.Range("K132").Value = 21606.1758
$ws.Range("L132").Value = 27591
$ws.Range("M132").Value = -19076.1758
$ws.Range("N132").Value = -32651

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 17927.477
$ws.Range("J54").Value = 21748.084
$ws.Range("L54").Value = 21748.084
$ws.Range("N54").Value = -22788.084
$ws.Range("H126").Value = 5266.1035
$ws.Range("I126").Value = 5702.577
$ws.Range("J126").Value = 1483.3334
$ws.Range("K126").Value = 17107.731
$ws.Range("L126").Value = 4450.0002
$ws.Range("M126").Value = -14637.731
$ws.Range("N126").Value = -9390.0002
$ws.Range("H132").Value = 18078.684
$ws.Range("I132").Value = 9544.187
$ws.Range("J132").Value = 39665.94
$ws.Range("K132").Value = 28632.561
$ws.Range("L132").Value = 118997.82
$ws.Range("M132").Value = -26102.561
$ws.Range("N132").Value = -124057.82
$ws.Range("H136").Value = 2107.6562
$ws.Range("I136").Value = 1460.8096
$ws.Range("J136").Value = 3342.5454
$ws.Range("K136").Value = 4382.4288
$ws.Range("L136").Value = 10027.6362
$ws.Range("M136").Value = -1832.4288
$ws.Range("N136").Value = -15127.6362
$ws.Range("H139").Value = 54215.312
$ws.Range("I139").Value = 45000
$ws.Range("J139").Value = 57287.082
$ws.Range("K139").Value = 45000
$ws.Range("L139").Value = 57287.082
$ws.Range("M139").Value = -39860
$ws.Range("N139").Value = -67567.08199999999
